$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): B11 4 -> 5, C11 -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total"): B12 112 -> 140, C12 0 -> -0, E12 "112/112" -> "140.0/140"
$ws.Range("B12").Value = 140
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "140.0/140"
